$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 512, shifting existing rows 512:575 down to 515:578.
$ws.Rows("512:514").Insert()

# Fill the 3 newly inserted rows (512, 513, 514) with the new data points.
# Columns: A Mercado ID, B Mercado, C Región, D Fecha, E Codreg, F Categoría ID,
# G Categoría, H Variedad, I Calidad, J Volumen, K Precio minimo, L Precio maximo,
# M Precio promedio ponderado, N Unidad de comercializacion, O Origen,
# P Precio $/Kg, Q Kg o Unidades, R Clasificacion

$ws.Range("A512:A514").Value = 11
$ws.Range("B512:B514").Value = "Vega Monumental Concepción"
$ws.Range("C512:C514").Value = "Bíobío"
$ws.Range("D512:D514").Value = 44918
$ws.Range("E512:E514").Value = 8
$ws.Range("F512:F514").Value = 100112020
$ws.Range("G512:G514").Value = "Tomate"
$ws.Range("H512:H514").Value = "Larga vida"
$ws.Range("N512:N514").Value = "$/bandeja 18 kilos"
$ws.Range("O512:O514").Value = "Región de O'Higgins"
$ws.Range("Q512:Q514").Value = 18
$ws.Range("R512:R514").Value = "Hortaliza"

$ws.Range("I512").Value = "Extra"
$ws.Range("J512").Value = 300
$ws.Range("K512").Value = 15000
$ws.Range("L512").Value = 15000
$ws.Range("M512").Value = 15000
$ws.Range("P512").Value = 833

$ws.Range("I513").Value = "Primera"
$ws.Range("J513").Value = 400
$ws.Range("K513").Value = 13000
$ws.Range("L513").Value = 13000
$ws.Range("M513").Value = 13000
$ws.Range("P513").Value = 722

$ws.Range("I514").Value = "Segunda"
$ws.Range("J514").Value = 300
$ws.Range("K514").Value = 11000
$ws.Range("L514").Value = 11000
$ws.Range("M514").Value = 11000
$ws.Range("P514").Value = 611
